$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows to append (Date serial, TotalUsers, LoggedInUsers, 0Err, 1Err, 2Err, 3-5Err, 6-10Err, >10Err)
$newRows = @(
    @(45993, 5620, 4351, 4064, 206, 43, 35, 3, 0),
    @(45994, 5617, 4388, 4026, 275, 45, 37, 5, 0)
)

$templateRow = 39
$startRow = 40

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    # Copy formatting/styles from the last existing data row down to the new row
    $ws.Range("A" + $templateRow + ":I" + $templateRow).Copy()
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122)

    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}

$excel.CutCopyMode = $false

# Update the selection to match the last new row, mirroring the original author's edit
$lastRow = $startRow + $newRows.Count - 1
$ws.Range("A" + $lastRow + ":I" + $lastRow).Select()
